$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.577462911605835
$ws.Range("B1").Value = 2.905165195465088
$ws.Range("C1").Value = 5.953004837036133
$ws.Range("D1").Value = 2.178303956985474
$ws.Range("E1").Value = 0.7835332751274109
